# Add 2022-Q4 data
# 1) Update the "总计" (Total) summary sheet: insert a new top data row for
#    2022-Q4 and shift the existing quarters down by one row.
# 2) Create a new "2022-Q4" worksheet (copied from "2022-Q3" so that it
#    inherits formatting/styles), positioned right after "总计", and fill it
#    in with the new quarter's fund holdings data.

$wb = $excel.ActiveWorkbook

# --- 1. Update 总计 sheet ---------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 8.72

$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 5
$total.Range("D3").Value = 10.13

$total.Range("B4").Value = "2022-Q2"
$total.Range("C4").Value = 4
$total.Range("D4").Value = 9.8

$total.Range("B5").Value = "2022-Q1"
$total.Range("C5").Value = 6
$total.Range("D5").Value = 11.66

$total.Range("B6").Value = "2021-Q4"
$total.Range("C6").Value = 6
$total.Range("D6").Value = 11.03

$total.Range("B7").Value = "2021-Q3"
$total.Range("C7").Value = 8
$total.Range("D7").Value = 16.41

$total.Range("B8").Value = "2021-Q2"
$total.Range("C8").Value = 6
$total.Range("D8").Value = 16.76

$total.Range("A9").Value = 7
$total.Range("B9").Value = "2021-Q1"
$total.Range("C9").Value = 8
$total.Range("D9").Value = 9.51

# --- 2. Create the new 2022-Q4 sheet -----------------------------------
$srcQ3 = $wb.Worksheets.Item("2022-Q3")
$srcQ3.Copy($srcQ3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# Fund rows keep the same fund codes/names as 2022-Q3; only the numeric
# columns (D..H) change.
# Row 2 - 010662 / 富国均衡优选混合
$q4.Range("D2").NumberFormat = "@"
$q4.Range("D2").Value = "42.60"
$q4.Range("E2").NumberFormat = "@"
$q4.Range("E2").Value = "89.38"
$q4.Range("F2").NumberFormat = "@"
$q4.Range("F2").Value = "8.05"
$q4.Range("G2").NumberFormat = "@"
$q4.Range("G2").Value = "3.4293"
$q4.Range("H2").Value = 1

# Row 3 - 010966 / 富国成长领航混合
$q4.Range("D3").NumberFormat = "@"
$q4.Range("D3").Value = "44.66"
$q4.Range("E3").NumberFormat = "@"
$q4.Range("E3").Value = "90.56"
$q4.Range("F3").NumberFormat = "@"
$q4.Range("F3").Value = "7.37"
$q4.Range("G3").NumberFormat = "@"
$q4.Range("G3").Value = "3.2914"
$q4.Range("H3").Value = 2

# Row 4 - 001985 / 富国低碳新经济混合A
$q4.Range("D4").NumberFormat = "@"
$q4.Range("D4").Value = "18.93"
$q4.Range("E4").NumberFormat = "@"
$q4.Range("E4").Value = "93.83"
$q4.Range("F4").NumberFormat = "@"
$q4.Range("F4").Value = "8.16"
$q4.Range("G4").NumberFormat = "@"
$q4.Range("G4").Value = "1.5447"
$q4.Range("H4").Value = 2

# Row 5 - 009693 / 富国积极成长一年定期开放混合
$q4.Range("D5").NumberFormat = "@"
$q4.Range("D5").Value = "12.47"
$q4.Range("E5").NumberFormat = "@"
$q4.Range("E5").Value = "98.05"
$q4.Range("F5").NumberFormat = "@"
$q4.Range("F5").Value = "2.46"
$q4.Range("G5").NumberFormat = "@"
$q4.Range("G5").Value = "0.3068"
$q4.Range("H5").Value = 9

# Row 6 - 011306 / 富国低碳新经济混合C
$q4.Range("D6").NumberFormat = "@"
$q4.Range("D6").Value = "1.87"
$q4.Range("E6").NumberFormat = "@"
$q4.Range("E6").Value = "93.83"
$q4.Range("F6").NumberFormat = "@"
$q4.Range("F6").Value = "8.16"
$q4.Range("G6").NumberFormat = "@"
$q4.Range("G6").Value = "0.1526"
$q4.Range("H6").Value = 2

$total.Select()
$total.Range("A1").Select()
